$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 18 with FizzBuzz entry
$ws.Range("A18").Value = "12th Aug"
$ws.Range("B18").Value = 412
$ws.Range("C18").Value = "FizzBuzz"
$ws.Range("D18").Value = "easy"
$ws.Range("E18").Value = "completed"

# Update selection to A19 as in the diff
$ws.Range("A19").Select()
